$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I and J (copy formatting from H1, same header style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in data rows 2-19: I column = 1, J column = same as H column
for ($r = 2; $r -le 19; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# Row 20 is special: I20 = 9, J20 = 9 (H20 stays 1)
$ws.Cells.Item(20, 9).Value = 9
$ws.Cells.Item(20, 10).Value = 9
